# Commit: "Modify the functionalities in page ShareSkill"
#
# - Add a new "CreditAmount" column (P) to the ShareSkill sheet, shifting the
#   previous "Active" header/value into a new column (Q).
# - Update the sample row's Startdate/Enddate to more recent dates.
# - Populate the new CreditAmount (P2) with a numeric value and the moved
#   Active value (Q2) with "Hidden".
# - Make ShareSkill the active/selected sheet (it was previously Profile).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# Row 1 headers: P1 becomes "CreditAmount", Q1 becomes the (now shifted) "Active" header
$ws.Range("P1").Value = "CreditAmount"
$ws.Range("Q1").Value = "Active"

# New column needs to fit its header text, like the other (bestFit) columns
$ws.Columns.Item(16).AutoFit()

# Row 2 data: refreshed Startdate / Enddate
$ws.Range("H2").Value = 44298
$ws.Range("I2").Value = 44328

# Row 2 data: CreditAmount value and the shifted Active ("Hidden") value
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = "Hidden"

# Make ShareSkill the active sheet/tab, with P5 selected
$ws.Activate()
$ws.Range("P5").Select() | Out-Null
